$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), O (Origen), P (Precio $/Kg)
# These correspond to a reshuffle of the weekly records already present in
# rows 2-14 (same underlying dataset, different week-to-row assignment).

$rows = @{
    2  = @{ D = 44418; J = 30; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí"; P = 500 }
    3  = @{ D = 44432; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí";  P = 467 }
    4  = @{ D = 44474; J = 45; K = 10000; L = 10000; M = 10000; O = "Provincia de Limarí"; P = 333 }
    5  = @{ D = 44376; J = 25; K = 18000; L = 18000; M = 18000; O = "Provincia de Limarí"; P = 600 }
    6  = @{ D = 44446; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí"; P = 467 }
    7  = @{ D = 44425; J = 35; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí"; P = 467 }
    8  = @{ D = 44449; J = 45; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 400 }
    9  = @{ D = 44435; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí"; P = 467 }
    10 = @{ D = 44435; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí";  P = 467 }
    11 = @{ D = 44421; J = 25; K = 15000; L = 16000; M = 15400; O = "Provincia de Limarí"; P = 513 }
    12 = @{ D = 44453; J = 50; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 400 }
    13 = @{ D = 44460; J = 45; K = 13000; L = 13000; M = 13000; O = "Provincia de Limarí"; P = 433 }
    14 = @{ D = 44467; J = 35; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí"; P = 400 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
}
